$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison (numeric MyForecast column updates) ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws1.Range("D4").Value = 103
$ws1.Range("D5").Value = 115
$ws1.Range("D10").Value = 107
$ws1.Range("D11").Value = 90
$ws1.Range("D12").Value = 99

# --- Sheet: Summary (these cells hold numeric/date-looking text, so force
#     the cell format to Text first to keep them stored as strings, matching
#     the source data which stores them as inline/shared text, not numbers
#     or dates) ---
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "2028"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "1257"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "772"

$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "90"

$ws2.Range("B15").NumberFormat = "@"
$ws2.Range("B15").Value = "2025-05-11"
